# Sync automático del tracker - append new prediction rows (74-78)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Predictions")

$rows = @(
    @{ A="2025-09-07"; B="Major League Soccer";  C="Houston Dynamo";      D="Los Angeles Galaxy";     E="Home Win"; F="69.42%"; G=1.85; H="27.14%"; I=2;   J=0.033443613920573;  K=0.33443613920573;  L="Pending" },
    @{ A="2025-09-07"; B="Major League Soccer";  C="Chicago Fire";        D="New England Revolution"; E="Home Win"; F="72.83%"; G=1.75; H="26.19%"; I=2.2; J=0.03661302357767575; K=0.3661302357767575; L="Pending" },
    @{ A="2025-09-07"; B="Major League Soccer";  C="St. Louis City";      D="FC Dallas";              E="Home Win"; F="57.45%"; G=2.05; H="16.60%"; I=1;   J=0.01692920300529338; K=0.1692920300529338; L="Pending" },
    @{ A="2025-09-07"; B="Liga de Expansión MX"; C="CDS Tampico Madero";  D="Tapatío";                E="Home Win"; F="77.01%"; G=1.62; H="23.51%"; I=2.4; J=0.03993199911892363; K=0.3993199911892363; L="Pending" },
    @{ A="2025-09-07"; B="Liga de Expansión MX"; C="Tepatitlán";          D="Leones Negros UDG";      E="Home Win"; F="57.80%"; G=2.1;  H="20.16%"; I=1.1; J=0.01942729720001626; K=0.1942729720001626; L="Pending" }
)

$startRow = 74
# Columns that must stay literal text even though their content looks
# like a date/percentage number (Excel would otherwise auto-convert them).
$textCols = @(1, 6, 8)   # A, F, H
$blankCols = @(13, 14, 15, 16, 17)  # M, N, O, P, Q

$rowIndex = $startRow
foreach ($row in $rows) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L")
    for ($c = 1; $c -le 12; $c++) {
        $colLetter = $cols[$c - 1]
        $value = $row[$colLetter]
        $cell = $ws.Cells.Item($rowIndex, $c)
        if ($textCols -contains $c) {
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }

    foreach ($c in $blankCols) {
        $cell = $ws.Cells.Item($rowIndex, $c)
        $cell.NumberFormat = "General"
        $cell.Style = "Normal"
    }

    $rowIndex = $rowIndex + 1
}
